$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.813.18'
$ws.Range('E2').Value = '  +0.88%  '
$ws.Range('D3').Value = '3.311.67'
$ws.Range('E3').Value = '  +4.98%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '603.71'
$ws.Range('E5').Value = '  +1.69%  '
$ws.Range('D6').Value = '142.41'
$ws.Range('E6').Value = '  +2.55%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '3.309.49'
$ws.Range('E8').Value = '  +5.10%  '
$ws.Range('D9').Value = '0.519'
$ws.Range('E9').Value = '  +0.46%  '
$ws.Range('E10').Value = '  +2.16%  '
$ws.Range('D11').Value = '5.47'
$ws.Range('E11').Value = '  +3.46%  '
$ws.Range('E12').Value = '  +2.07%  '
$ws.Range('E13').Value = '  +0.80%  '
$ws.Range('D14').Value = '34.70'
$ws.Range('E14').Value = '  +1.15%  '
$ws.Range('D15').Value = '3.859.91'
$ws.Range('E15').Value = '  +5.16%  '
$ws.Range('E16').Value = '  +0.07%  '
$ws.Range('D17').Value = '3.311.65'
$ws.Range('E17').Value = '  +5.23%  '
$ws.Range('D18').Value = '63.884.28'
$ws.Range('E18').Value = '  +1.07%  '
$ws.Range('E19').Value = '  +2.46%  '
$ws.Range('D20').Value = '480.44'
$ws.Range('E20').Value = '  +1.08%  '
$ws.Range('D21').Value = '14.18'
$ws.Range('E21').Value = '  +0.19%  '
$ws.Range('D22').Value = '0.732'
$ws.Range('E22').Value = '  +4.50%  '
$ws.Range('D23').Value = '8.06'
$ws.Range('E23').Value = '  +4.59%  '
$ws.Range('D24').Value = '13.75'
$ws.Range('E24').Value = '  +5.61%  '
$ws.Range('D25').Value = '84.56'
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('E27').Value = '  +1.56%  '
$ws.Range('D28').Value = '7.35'
$ws.Range('E28').Value = '  +4.35%  '
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('D30').Value = '8.12'
$ws.Range('E30').Value = '  +0.92%  '
$ws.Range('E31').Value = '  +2.14%  '
$ws.Range('D32').Value = '28.77'
$ws.Range('E32').Value = '  +6.77%  '
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('E34').Value = '  +0.39%  '
$ws.Range('E35').Value = '  +3.68%  '
$ws.Range('E36').Value = '  +3.13%  '
$ws.Range('D37').Value = '53.35'
$ws.Range('E37').Value = '  +1.71%  '
$ws.Range('D38').Value = '0.0₃0741'
$ws.Range('E38').Value = '  +5.18%  '
$ws.Range('E39').Value = '  +2.68%  '
$ws.Range('D40').Value = '434.12'
$ws.Range('E40').Value = '  +2.79%  '
$ws.Range('D41').Value = '3.068.74'
$ws.Range('E41').Value = '  +3.72%  '
$ws.Range('E42').Value = '  -0.44%  '
$ws.Range('D43').Value = '8.34'
$ws.Range('E43').Value = '  +0.99%  '
$ws.Range('E44').Value = '  +1.37%  '
$ws.Range('E45').Value = '  +0.36%  '
$ws.Range('E46').Value = '  +3.22%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = '26.37'
$ws.Range('E47').Value = '  +3.34%  '
$ws.Range('B48').Value = 'Arweave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D48').Value = '36.25'
$ws.Range('E48').Value = '  +12.54%  '
$ws.Range('D50').Value = '126.34'
$ws.Range('E50').Value = '  +5.13%  '
$ws.Range('E51').Value = '  +0.85%  '
